$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order row (row 7) - values are stored as plain text in the source
# file (same as every other cell in the sheet), so the numeric-looking
# ones need an explicit text/apostrophe marker or Excel will reinterpret
# them as numbers/currency.
$ws.Range("A7").Value = 'PS 1077861'
$ws.Range("B7").Value = 'Wrapping Sheets - 15x15 (poly)'
$ws.Range("C7").Value = "'1"
$ws.Range("D7").Value = '''$41.40'
$ws.Range("E7").Value = '''$41.40'
